# Adds a new "LAT testing result of MNIST" block to the sheet:
#  - renames the header of column C ("逐层对抗训练模型" -> "原始逐层对抗训练模型")
#  - adds two new result columns H, I for the new LAT variants with their
#    header labels (row1/row2) and accuracy numbers (rows 5-8)
#  - moves the current selection to E24

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Give the two new columns (H, I) the same look as the neighbouring
#     header cells (bold style "s=3") before filling them in -------------
$ws.Range("G1:G2").Copy()
$ws.Range("H1:H2").PasteSpecial(-4122)
$ws.Range("I1:I2").PasteSpecial(-4122)
$ws.Range("J1:J2").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Row 1: model-family headers ------------------------------------------
$ws.Range("H1").Value = "LAT_model"
$ws.Range("I1").Value = "LAT_model"

# --- Row 2: run/column headers --------------------------------------------
$ws.Range("H2").Value = "逐层对抗训练模型(change the bp method,without ZeroReg)"
$ws.Range("I2").Value = "逐层对抗训练模型(change the bp method,ZeroReg every steps)"

# --- Column C header text changes (same column, new label) ---------------
$ws.Range("C2").Value = "原始逐层对抗训练模型"

# --- Accuracy numbers for the new columns (rows 5-8) ----------------------
$ws.Range("H5").Value = 0.99
$ws.Range("I5").Value = 0.99

$ws.Range("H6").Value = 0.97
$ws.Range("I6").Value = 0.97

$ws.Range("H7").Value = 0.93
$ws.Range("I7").Value = 0.93

$ws.Range("H8").Value = 0.76
$ws.Range("I8").Value = 0.77

# --- Move the active selection, matching the saved view -------------------
[void]$ws.Range("E24").Select()
